$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -9
$ws.Range("F3").Value = -3
$ws.Range("F4").Value = -3
$ws.Range("F8").Value = -8
$ws.Range("F9").Value = -9
$ws.Range("F14").Value = 3
